$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Insert programban"
$ws.Range("B11").Value = "Success"
$ws.Range("A12").Value = "Select programban"
$ws.Range("B12").Value = "Success"

$ws.Range("A13").Select()
